$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1 / sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6497
$wsExhibit.Range("F10").Value = 83
$wsExhibit.Range("F13").Value = 380
$wsExhibit.Range("F14").Value = 951
$wsExhibit.Range("F15").Value = 3195
$wsExhibit.Range("F19").Value = 25

# Sheet "全部类型" (sheetId 4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6497
$wsAll.Range("F11").Value = 83
$wsAll.Range("F14").Value = 380
$wsAll.Range("F15").Value = 951
$wsAll.Range("F16").Value = 3195
$wsAll.Range("F20").Value = 25
